{"js": "// Fix the typo \"collaboratio\" -> \"collaboration\" and remove the two\n// bullet paragraphs about Linux cluster environments and about storing/\n// transferring large-scale genomics data (per commit message \"updating\n// to remove linux\"), while keeping the \"Managed the installation and\n// configuration of R packages...\" bullet in between them.\n\nconst body = context.document.body;\n\n// 1) Fix the typo \"collaboratio\" -> \"collaboration\".\nconst typoResults = body.search(\"collaboratio,\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"collaboration,\", \"Replace\");\n}\n\n// 2) Remove the \"Developed and maintained Linux cluster environments...\" bullet.\nconst linuxResults = body.search(\n  \"Developed and maintained Linux cluster environments, including the installation and management of bioinformatics tools for large-scale genomic analysis (Docker).\",\n  { matchCase: true }\n);\nlinuxResults.load(\"text\");\nawait context.sync();\n\nif (linuxResults.items.length > 0) {\n  const para = linuxResults.items[0].paragraphs.getFirst();\n  para.delete();\n}\n\n// 3) Remove the \"Stored, indexed, and transferred large-scale genomics...\" bullet.\nconst storedResults = body.search(\n  \"Stored, indexed, and transferred large-scale genomics and clinical trial data, utilizing high-performance storage solutions and cloud platforms (AWS).\",\n  { matchCase: true }\n);\nstoredResults.load(\"text\");\nawait context.sync();\n\nif (storedResults.items.length > 0) {\n  const para2 = storedResults.items[0].paragraphs.getFirst();\n  para2.delete();\n}\n\nawait context.sync();\n", "ps1": "# Fix the typo \"collaboratio\" -> \"collaboration\" and remove the two\n# bullet paragraphs about Linux cluster environments and about storing/\n# transferring large-scale genomics data (per commit message \"updating\n# to remove linux\"), while keeping the \"Managed the installation and\n# configuration of R packages...\" bullet in between them.\n\n$d = $word.ActiveDocument\n\n# 1) Fix the typo \"collaboratio\" -> \"collaboration\".\n$find = $d.Content.Find\n$find.Execute(\"collaboratio,\", $false, $false, $false, $false, $false, $true, 1, $false, \"collaboration,\", 2) | Out-Null\n\n# 2) Remove the two obsolete bullet paragraphs about Linux / data storage,\n#    leaving the \"Managed the installation...\" bullet intact.\n$targets = @(\n    \"Developed and maintained Linux cluster environments, including the installation and management of bioinformatics tools for large-scale genomic analysis (Docker).\",\n    \"Stored, indexed, and transferred large-scale genomics and clinical trial data, utilizing high-performance storage solutions and cloud platforms (AWS).\"\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n        if ($text -eq $target) {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n"}
